$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting of row 25 (a fully-populated "data" row) onto row 26,
# so row 26 gets the same styling Excel applied to the new log entry row.
$ws.Range("A25:F25").Copy()
$ws.Range("A26:F26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new log entry (row 26) ---
$ws.Range("A26").Value = "Fixed tag issues"
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = "5/19/2025"
$ws.Range("D26").Value = "I started off by switching out the tag_id params for tagName params, since it would not be possible to create tags when just passing tag ID's. So now you just pass in tag names, if the tag already exists, it will add that tag, and if it doesnt it will add the tag to the tags in the DB and add it to the facility. I also added validation checks for new tags to check if a string is given and to check if its not an empty string. I aslo added a validateFacilityId method to not write duplicate code. This all took me longer then expected but oh well. The only thing left, besides the bonus stuff, is the use of the model thingies. Im going to try my best to understand how to implement this. im going to continue tomorrow since its getting late haha."

# --- Row height changes ---
$ws.Range("A9").RowHeight = 25.5
$ws.Range("A26").RowHeight = 65.25

# --- View state update: scroll down and move selection ---
$ws.Range("D35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
